# Update Name of Algo
# Applies updated KNN imputation results to specific cells in Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "C3"   = -13.406
    "B7"   = 5.697000000000001
    "A10"  = -21.604
    "A12"  = -21.515
    "B15"  = 5.072000000000001
    "A18"  = -21.691
    "C18"  = -10.732
    "C19"  = -11.706
    "B20"  = 7.217000000000001
    "C27"  = -13.022
    "B29"  = 5.645
    "B30"  = 6.02
    "B31"  = 6.367
    "A37"  = -20.287
    "B40"  = 8.870999999999999
    "C42"  = -12.533
    "C44"  = -12.324
    "C47"  = -12.357
    "A55"  = -21.875
    "C58"  = -12.621
    "A68"  = -21.581
    "B68"  = 5.242
    "C73"  = -12.57
    "B76"  = 6.654000000000001
    "A77"  = -20.609
    "A78"  = -19.842
    "B87"  = 5.429
    "B88"  = 5.75
    "C95"  = -11.598
    "B96"  = 6.484
    "B98"  = 5.358
    "B101" = 9.096
    "C101" = -12.234
    "B102" = 7.739999999999999
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
